$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 109, pushing the existing rows 109-129 down to 110-130.
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with a fresh price observation
# (same product/quality/origin as the last row, but a new date & prices).
$ws.Range("A109").Value = 4
$ws.Range("B109").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C109").Value = "Los Lagos"
$ws.Range("D109").Value = 45009
$ws.Range("D109").NumberFormat = $ws.Range("D110").NumberFormat
$ws.Range("E109").Value = 10
$ws.Range("F109").Value = "Fruta"
$ws.Range("G109").Value = 100104
$ws.Range("H109").Value = "Frutos de pepita"
$ws.Range("I109").Value = 100104003
$ws.Range("J109").Value = "Membrillo"
$ws.Range("K109").Value = "Champion"
$ws.Range("L109").Value = "Primera"
$ws.Range("M109").Value = 400
$ws.Range("N109").Value = 17000
$ws.Range("O109").Value = 18000
$ws.Range("P109").Value = 17500
$ws.Range("Q109").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R109").Value = "Región de O'Higgins"
$ws.Range("S109").Value = 972
$ws.Range("T109").Value = 18
